# Daily attendance processing - 2026-01-05 18:45:04
# Reorder the "Recorded By" (column G) contributor list for rows whose
# value matches the known patterns, e.g.
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   "System, system, backup@backdoor.com" -> "System, backup@backdoor.com, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    } elseif ($v -eq "System, system, backup@backdoor.com") {
        $cell.Value = "System, backup@backdoor.com, system"
    }
}
